$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '29.483.37'
$ws.Range("E2").Value = '  +0.40%  '
$ws.Range("D3").Value = '1.848.00'
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.9980'
$ws.Range("E4").Value = '  -0.17%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '240.45'
$ws.Range("E5").Value = '  +0.10%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.6294'
$ws.Range("E6").Value = '  +0.18%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.9995'
$ws.Range("E7").Value = '  -0.08%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.07481'
$ws.Range("E8").Value = '  -1.60%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.2911'
$ws.Range("E9").Value = '  -0.03%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '24.63'
$ws.Range("E10").Value = '  +0.26%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07741'
$ws.Range("D12").Value = '1.847.16'
$ws.Range("E12").Value = '  -0.05%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '5.015'
$ws.Range("E13").Value = '  +0.01%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.6816'
$ws.Range("E14").Value = '  +0.28%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.00001044'
$ws.Range("E15").Value = '  +0.01%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '82.19'
$ws.Range("E16").Value = '  -1.04%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '6.250'
$ws.Range("E17").Value = '  +2.23%  '
$ws.Range("D18").Value = '29.479.82'
$ws.Range("E18").Value = '  +0.35%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '229.43'
$ws.Range("E19").Value = '  +0.24%  '
$ws.Range("E20").Value = '  +0.59%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '0.9993'
$ws.Range("E21").Value = '  -0.05%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '7.538'
$ws.Range("E22").Value = '  +1.29%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.9998'
$ws.Range("E23").Value = '  +0.03%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '159.39'
$ws.Range("E24").Value = '  +0.30%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '8.523'
$ws.Range("E25").Value = '  +1.02%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.1369'
$ws.Range("E26").Value = '  -1.56%  '
$ws.Range("E27").Value = '  -0.71%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '0.06524'
$ws.Range("E28").Value = '  +16.13%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.422'
$ws.Range("E29").Value = '  -0.71%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.488'
$ws.Range("E30").Value = '  +1.27%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '4.103'
$ws.Range("E31").Value = '  -0.13%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '4.105'
$ws.Range("E32").Value = '  +1.55%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.836'
$ws.Range("E33").Value = '  +0.59%  '
$ws.Range("E34").Value = '  -0.93%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.7000'
$ws.Range("E35").Value = '  +0.55%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '2.577'
$ws.Range("E36").Value = '  -0.18%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.01865'
$ws.Range("E37").Value = '  +1.76%  '
$ws.Range("D38").Value = '1.266.17'
$ws.Range("E38").Value = '  +2.55%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '2.840'
$ws.Range("E39").Value = '  +4.35%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '6.823'
$ws.Range("E40").Value = '  +6.57%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.9364'
$ws.Range("E41").Value = '  +4.04%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.9998'
$ws.Range("E42").Value = '  +0.01%  '
$ws.Range("D43").Value = '2.020.20'
$ws.Range("E43").Value = '  +0.70%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '101.35'
$ws.Range("E44").Value = '  +0.02%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '66.24'
$ws.Range("E45").Value = '  +1.24%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '1.736'
$ws.Range("E46").Value = '  +3.35%  '
$ws.Range("B47").Value = 'BabyDogeCoin'
$ws.Range("C47").Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.00000000118'
$ws.Range("E47").Value = '  +2.95%  '
$ws.Range("B48").Value = 'Aptos'
$ws.Range("C48").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '7.101'
$ws.Range("E48").Value = '  -0.51%  '
$ws.Range("B49").Value = 'Algorand'
$ws.Range("C49").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.1166'
$ws.Range("E49").Value = '  +1.35%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '9.017'
$ws.Range("E50").Value = '  +0.40%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.3961'
$ws.Range("E51").Value = '  -0.83%  '
